$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell K1 - copy header style/format from J1 (bold, border, centered), then set its text
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("K1").Value = "intervention_type"

# New column K: intervention_type per clinical trial row.
# Rows with no known intervention type are intentionally left blank.
$ws.Range("K2").Value = "PROCEDURE"
$ws.Range("K3").Value = "BIOLOGICAL"
$ws.Range("K4").Value = "OTHER"
$ws.Range("K5").Value = "OTHER"
$ws.Range("K6").Value = "DRUG"
$ws.Range("K7").Value = "OTHER"
$ws.Range("K8").Value = "DRUG"
$ws.Range("K9").Value = "DRUG"
$ws.Range("K10").Value = "DEVICE"
$ws.Range("K11").Value = "OTHER"
$ws.Range("K12").Value = "BEHAVIORAL"
$ws.Range("K13").Value = "OTHER"
$ws.Range("K14").Value = "DRUG"
$ws.Range("K15").Value = "PROCEDURE"
$ws.Range("K16").Value = "DRUG"
$ws.Range("K17").Value = "OTHER"
$ws.Range("K18").Value = "DEVICE"
$ws.Range("K19").Value = "PROCEDURE"
$ws.Range("K20").Value = "OTHER"
$ws.Range("K21").Value = "OTHER"
$ws.Range("K22").Value = "OTHER"
$ws.Range("K23").Value = "RADIATION"
$ws.Range("K24").Value = "PROCEDURE"
$ws.Range("K25").Value = "PROCEDURE"
$ws.Range("K26").Value = "DRUG"
$ws.Range("K27").Value = "BIOLOGICAL"
$ws.Range("K28").Value = "DEVICE"
$ws.Range("K29").Value = "OTHER"
$ws.Range("K30").Value = "OTHER"
$ws.Range("K31").Value = "OTHER"
$ws.Range("K32").Value = "BEHAVIORAL"
$ws.Range("K33").Value = "GENETIC"
$ws.Range("K35").Value = "DRUG"
$ws.Range("K36").Value = "OTHER"
$ws.Range("K37").Value = "DRUG"
$ws.Range("K38").Value = "DRUG"
$ws.Range("K39").Value = "DRUG"
$ws.Range("K40").Value = "DEVICE"
$ws.Range("K41").Value = "DIAGNOSTIC_TEST"
$ws.Range("K42").Value = "PROCEDURE"
$ws.Range("K43").Value = "OTHER"
$ws.Range("K44").Value = "DRUG"
$ws.Range("K45").Value = "BEHAVIORAL"
$ws.Range("K46").Value = "PROCEDURE"
$ws.Range("K47").Value = "DRUG"
$ws.Range("K48").Value = "DRUG"
$ws.Range("K49").Value = "BIOLOGICAL"
$ws.Range("K50").Value = "DRUG"
$ws.Range("K51").Value = "OTHER"
$ws.Range("K52").Value = "OTHER"
$ws.Range("K53").Value = "OTHER"
$ws.Range("K54").Value = "OTHER"
$ws.Range("K55").Value = "DRUG"
$ws.Range("K56").Value = "DEVICE"
$ws.Range("K57").Value = "OTHER"
$ws.Range("K58").Value = "PROCEDURE"
$ws.Range("K59").Value = "DRUG"
$ws.Range("K60").Value = "OTHER"
$ws.Range("K61").Value = "OTHER"
$ws.Range("K62").Value = "DRUG"
$ws.Range("K63").Value = "DIAGNOSTIC_TEST"
$ws.Range("K64").Value = "PROCEDURE"
$ws.Range("K65").Value = "DEVICE"
$ws.Range("K66").Value = "OTHER"
$ws.Range("K67").Value = "PROCEDURE"
$ws.Range("K68").Value = "OTHER"
$ws.Range("K69").Value = "DEVICE"
$ws.Range("K70").Value = "DEVICE"
$ws.Range("K71").Value = "DEVICE"
$ws.Range("K72").Value = "OTHER"
$ws.Range("K73").Value = "DRUG"
$ws.Range("K74").Value = "OTHER"
$ws.Range("K75").Value = "OTHER"
$ws.Range("K76").Value = "DRUG"
$ws.Range("K77").Value = "OTHER"
$ws.Range("K78").Value = "OTHER"
$ws.Range("K79").Value = "DRUG"
$ws.Range("K80").Value = "OTHER"
$ws.Range("K81").Value = "DRUG"
$ws.Range("K82").Value = "DEVICE"
$ws.Range("K85").Value = "DRUG"
$ws.Range("K86").Value = "PROCEDURE"
$ws.Range("K87").Value = "DRUG"
$ws.Range("K88").Value = "OTHER"
$ws.Range("K89").Value = "BIOLOGICAL"
$ws.Range("K90").Value = "BIOLOGICAL"
$ws.Range("K91").Value = "OTHER"
$ws.Range("K92").Value = "BIOLOGICAL"
$ws.Range("K93").Value = "DRUG"

Write-Host "Added intervention_type column (K) with values for all sponsors rows."
